$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.195697069168091
$ws.Range("B1").Value = 2.386084079742432
$ws.Range("C1").Value = 4.366960048675537
$ws.Range("D1").Value = 2.775581836700439
$ws.Range("E1").Value = 1.112748503684998
